$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 55: 2025/10/03, 金, 17, 33
# Force column A to be treated as text (not auto-converted to a date
# serial) to match the existing "2025/10/03"-style text entries above,
# then drop the temporary number-format override so no stray style is
# left attached to the cell.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2025/10/03"
$ws.Range("A55").ClearFormats()

$ws.Range("B55").Value = "金"
$ws.Range("C55").Value = 17
$ws.Range("D55").Value = 33
